$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Oct 09 22:47:33 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 22:47:46 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 22:47:58 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 09 22:48:10 EDT 2023"
